# Regenerate save_data to use K (strikeouts) instead of Strike# column,
# updating the recalculated K values (column G) for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 8
    4  = 5
    5  = 7
    6  = 8
    7  = 4
    8  = 8
    9  = 5
    10 = 3
    11 = 11
    12 = 2
    13 = 7
    14 = 5
    15 = 4
    16 = 6
    17 = 8
    18 = 4
    19 = 4
    20 = 6
    21 = 10
    22 = 1
    23 = 9
    24 = 8
    25 = 6
    26 = 7
    27 = 4
    28 = 6
    29 = 8
    30 = 8
    31 = 4
    32 = 7
    33 = 6
    34 = 8
    35 = 4
    36 = 4
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
